# Applies the "updated 4.0 files and mdl" revision:
#  - About!C1: bump the "last updated" date serial 45294 -> 45379 (2024-01-03 -> 2024-03-28)
#  - FPIEBP!B3:D3 (hard coal production/imports/exports priorities): 3,2,1 -> 1,3,2
#  - Move the live selection on FPIEBP from F4 to E3
#  - Scroll the About sheet so row 6 is pinned at the top (topLeftCell = A6)

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")

# --- About sheet: refresh the "last updated" date stamp in C1 ---
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet: re-prioritize "hard coal" (row 3) production/imports/exports ---
$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# --- Cosmetic: restore the view state captured with the edit ---
$wsAbout.Activate()
$wsAbout.Range("A6").Select()
$excel.ActiveWindow.ScrollRow = 6

$wsFPIEBP.Activate()
$wsFPIEBP.Range("E3").Select()
